$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 17: Inscritos (E) 89 -> 90
$ws.Range("E17").Value = 90

# Row 19: Inscritos (E) 37 -> 38
$ws.Range("E19").Value = 38

# Row 37: Inscritos (E) 39 -> 40, Pagos (F) 18 -> 19, Inscricoes homologadas (H) 18 -> 19
$ws.Range("E37").Value = 40
$ws.Range("F37").Value = 19
$ws.Range("H37").Value = 19

# Row 66: Pagos (F) 14 -> 15, Inscricoes homologadas (H) 14 -> 15
$ws.Range("F66").Value = 15
$ws.Range("H66").Value = 15

# Row 70: Inscritos (E) 29 -> 30
$ws.Range("E70").Value = 30

# Row 72: Inscritos (E) 29 -> 30
$ws.Range("E72").Value = 30

# Row 77: Inscritos (E) 43 -> 44
$ws.Range("E77").Value = 44
